# Column category #7 - add a "Category" column (F) and rename/clean the
# existing "speedflow..." comments in column E, per commit:
#   "Column category #7 Added column category did some localization and
#    cleaning thus changed the example files"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Localize / clean up column E ("speedflowXXX" -> "Comment ...") ---
$ws.Range("E1").Value = "Comment ko"
$ws.Range("E2").Value = "Comment koko"
$ws.Range("E3").Value = "Comment kokoko"

# --- New column F: "Category N" per row ---
$ws.Range("F1").Value = "Category 1"
$ws.Range("F2").Value = "Category 2"
$ws.Range("F3").Value = "Category 3"

# Give the new column a sensible width (characters), matching the other
# data columns' sizing convention.
$ws.Columns.Item(6).ColumnWidth = 20.5

# Reflect the new column in the view: scroll right a bit and leave the
# selection on the newly added data (F2).
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("F2").Select()
